$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily log entry (2026/02/06, 金, 11, 60) is inserted as row 763,
# pushing the existing rows 763:804 down to 764:805 (dimension grows from
# A1:D804 to A1:D805).
$ws.Rows.Item(763).Insert()

# Column A holds dates as literal text (e.g. "2026/12/29"), not real Excel
# dates. Assigning a date-shaped string directly would get auto-converted
# to a date serial number, so prefix it with an apostrophe to force text
# entry, then reset the style to Normal so the cell doesn't end up with a
# stray quote-prefix style (matches the un-styled data cells elsewhere in
# the column).
$ws.Cells.Item(763, 1).Value = "'2026/02/06"
$ws.Cells.Item(763, 1).Style = "Normal"

$ws.Cells.Item(763, 2).Value = "金"
$ws.Cells.Item(763, 3).Value = 11
$ws.Cells.Item(763, 4).Value = 60
